# Commit: Comment => Annotation in XLSX
# The header cell H1 on Sheet1 changes its text from "Comment" to "Annotation".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H1").Value = "Annotation"
